# Updates cryptos list data (Price and Volume(1h) columns) per upstream scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.720.56"
Set-TextValue "E2" "  -0.47%  "
Set-TextValue "D3" "1.895.00"
Set-TextValue "E3" "  +1.44%  "
Set-TextValue "E4" "  -1.10%  "
Set-TextValue "D5" "313.22"
Set-TextValue "E5" "  +0.08%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  -0.90%  "
Set-TextValue "D7" "0.4853"
Set-TextValue "E7" "  +0.79%  "
Set-TextValue "D8" "0.3793"
Set-TextValue "E8" "  -0.66%  "
Set-TextValue "D9" "0.07337"
Set-TextValue "E9" "  -0.37%  "
Set-TextValue "D10" "0.9152"
Set-TextValue "E10" "  -2.38%  "
Set-TextValue "E11" "  -2.20%  "
Set-TextValue "D12" "0.07692"
Set-TextValue "E12" "  -1.39%  "
Set-TextValue "D13" "1.866.24"
Set-TextValue "E13" "  -0.66%  "
Set-TextValue "D14" "5.471"
Set-TextValue "E14" "  -0.17%  "
Set-TextValue "D15" "6.602"
Set-TextValue "E15" "  +0.06%  "
Set-TextValue "D16" "90.96"
Set-TextValue "E16" "  +0.33%  "
Set-TextValue "E17" "  -1.15%  "
Set-TextValue "D18" "0.000008811"
Set-TextValue "E18" "  -0.65%  "
Set-TextValue "E19" "  -0.79%  "
Set-TextValue "D20" "27.759.14"
Set-TextValue "E20" "  -1.07%  "
Set-TextValue "D21" "14.47"
Set-TextValue "E21" "  -2.21%  "
Set-TextValue "D22" "5.118"
Set-TextValue "E22" "  -0.06%  "
Set-TextValue "D23" "2.118.00"
Set-TextValue "E23" "  -0.69%  "
Set-TextValue "D24" "10.75"
Set-TextValue "E24" "  -0.57%  "
Set-TextValue "D25" "1.904"
Set-TextValue "E25" "  -1.79%  "
Set-TextValue "D26" "153.90"
Set-TextValue "E26" "  -1.55%  "
Set-TextValue "E27" "  -1.06%  "
Set-TextValue "D28" "2.137"
Set-TextValue "E28" "  +4.28%  "
Set-TextValue "E29" "  -0.04%  "
Set-TextValue "D30" "4.899"
Set-TextValue "E30" "  -1.42%  "
Set-TextValue "D31" "0.08911"
Set-TextValue "E31" "  +0.00%  "
Set-TextValue "D32" "3.154"
Set-TextValue "E32" "  -5.43%  "
Set-TextValue "D33" "1.224"
Set-TextValue "E33" "  +0.54%  "
Set-TextValue "D34" "0.7647"
Set-TextValue "E34" "  +0.18%  "
Set-TextValue "D35" "4.646"
Set-TextValue "E35" "  -0.21%  "
Set-TextValue "E36" "  +0.00%  "
Set-TextValue "D37" "2.524"
Set-TextValue "E37" "  -7.39%  "
Set-TextValue "D38" "1.092"
Set-TextValue "E38" "  -3.71%  "
Set-TextValue "D39" "0.05271"
Set-TextValue "E39" "  -1.87%  "
Set-TextValue "D40" "0.5468"
Set-TextValue "E40" "  -2.90%  "
Set-TextValue "E41" "  -0.17%  "
Set-TextValue "D42" "6.918"
Set-TextValue "D43" "8.455"
Set-TextValue "E43" "  -1.09%  "
Set-TextValue "E44" "  -0.85%  "
Set-TextValue "D45" "110.07"
Set-TextValue "E45" "  +4.91%  "
Set-TextValue "D46" "10.61"
Set-TextValue "E46" "  -1.53%  "
Set-TextValue "E47" "  -2.16%  "
Set-TextValue "D48" "1.001"
Set-TextValue "E48" "  -0.89%  "
Set-TextValue "D49" "1.636"
Set-TextValue "E49" "  -2.35%  "
Set-TextValue "D50" "67.31"
Set-TextValue "E50" "  -0.48%  "
